# The edit rotates the record data (columns A,B,D,E,F,G,H,I,J,P,Q,R) among
# rows 2-7 of the sheet: row 2's data -> row 4, row 4's data -> row 5,
# row 5's data -> row 7, row 7's data -> row 3, row 3's data -> row 6,
# row 6's data -> row 2 (a single 6-cycle). All other columns/rows/styles
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "I", "J", "P", "Q", "R")

# Capture the current ("before") values of every touched column for rows 2-7
# so that later writes don't clobber data we still need to read.
$data = @{}
foreach ($r in 2..7) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $data[$r] = $rowVals
}

# after_row -> before_row (where the data placed into after_row comes from)
$mapping = @{
    2 = 6
    3 = 7
    4 = 2
    5 = 4
    6 = 3
    7 = 5
}

foreach ($destRow in 2..7) {
    $srcRow = $mapping[$destRow]
    $srcVals = $data[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
